# Update AIPE portfolio data workbook:
#  - Sheet 1 (大智投资组合): refresh "修改时间" (E) timestamps, tweak D5 allocation
#  - Sheet 2 (大成投资组合): drop the "600580 卧龙电驱" holding (row 8), refresh
#    remaining allocations (D) and timestamps (E)
#  - Sheet 3 (我的投资组合): refresh "修改时间" (G) timestamps, tweak a few F
#    allocation values
#
# New timestamp stamped on every touched row, matching the commit's edit time.
$newTimestamp = "202509250137"

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: 大智投资组合
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

$ws1.Range("D5").Value = 5.27

$sheet1Rows = 2..9
foreach ($r in $sheet1Rows) {
    $cell = $ws1.Cells.Item($r, 5)
    # Force text storage (the stamp is a text column, not a number) the same
    # way typing into a pre-formatted "Text" cell would, then drop the
    # quote-prefix formatting remnant so the cell stays plain/default-styled.
    $cell.NumberFormat = "@"
    $cell.Value = $newTimestamp
    $cell.ClearFormats()
}

# ---------------------------------------------------------------------------
# Sheet 2: 大成投资组合 - remove the 卧龙电驱 (600580) holding entirely
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)

$ws2.Rows.Item(8).Delete()

# After the delete, rows 9-12 shifted up to 8-11; their B/C contents already
# line up with the new layout, only the allocation % (D) needs refreshing.
$ws2.Range("D2").Value = 4.92
$ws2.Range("D4").Value = 4.74
$ws2.Range("D6").Value = 5.2
$ws2.Range("D8").Value = 10
$ws2.Range("D9").Value = 4.84
$ws2.Range("D11").Value = 1.01

$sheet2Rows = 2..11
foreach ($r in $sheet2Rows) {
    $cell = $ws2.Cells.Item($r, 5)
    $cell.NumberFormat = "@"
    $cell.Value = $newTimestamp
    $cell.ClearFormats()
}

# ---------------------------------------------------------------------------
# Sheet 3: 我的投资组合
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item(3)

$ws3.Range("F4").Value = 5.1
$ws3.Range("F8").Value = 5.39
$ws3.Range("F9").Value = 3.19

$sheet3Rows = 2..13
foreach ($r in $sheet3Rows) {
    $cell = $ws3.Cells.Item($r, 7)
    $cell.NumberFormat = "@"
    $cell.Value = $newTimestamp
    $cell.ClearFormats()
}

Write-Host "Portfolio data refreshed."
